$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# The sheet is an attendance log. Each "session" is a 4-row block:
#   row N   : date (col A) / time (col C)
#   row N+1 : "Sno" / "Name" / "Present" header
#   row N+2 : record 1 (Sno / Name / Present)
#   row N+3 : record 2 (Sno / Name / Present)
# followed by one blank separator row.
#
# Before this edit the most-recent block occupied rows 87-90 and its
# cells had no explicit style. Two new sessions were recorded after
# it (26/04/2019 02:48:37 and 26/04/2019 02:52:01), pushing the
# dimension out to row 100. Once a block is no longer the very last
# one written, it picks up the sheet's normal per-cell style - so we
# re-enter rows 87-90 (same values) to pick that up, then append the
# two new blocks.
# ------------------------------------------------------------------

# Re-write rows 87-90 so they stop being "freshly appended" and take
# on the same formatting as every earlier block in the sheet.
$ws.Range("A87:C90").ClearContents()
$ws.Range("A87").Value = "26/02/2019"
$ws.Range("C87").Value = "21:22:24"
$ws.Range("A88").Value = "Sno"
$ws.Range("B88").Value = "Name"
$ws.Range("C88").Value = "Present"
$ws.Range("A89").Value = 1
$ws.Range("B89").Value = "shamil"
$ws.Range("C89").Value = "no"
$ws.Range("A90").Value = 2
$ws.Range("B90").Value = "abhi"
$ws.Range("C90").Value = "no"

# New session recorded 26/04/2019 02:48:37 (rows 92-95, row 91 blank)
$ws.Range("A92").Value = "26/04/2019"
$ws.Range("C92").Value = "02:48:37"
$ws.Range("A93").Value = "Sno"
$ws.Range("B93").Value = "Name"
$ws.Range("C93").Value = "Present"
$ws.Range("A94").Value = 1
$ws.Range("B94").Value = "abhi"
$ws.Range("C94").Value = "yes"
$ws.Range("A95").Value = 2
$ws.Range("B95").Value = "shamil"
$ws.Range("C95").Value = "no"

# Newest session recorded 26/04/2019 02:52:01 (rows 97-100, row 96 blank)
$ws.Range("A97").Value = "26/04/2019"
$ws.Range("C97").Value = "02:52:01"
$ws.Range("A98").Value = "Sno"
$ws.Range("B98").Value = "Name"
$ws.Range("C98").Value = "Present"
$ws.Range("A99").Value = 1
$ws.Range("B99").Value = "abhi"
$ws.Range("C99").Value = "no"
$ws.Range("A100").Value = 2
$ws.Range("B100").Value = "shamil"
$ws.Range("C100").Value = "no"
